# Auto-generated Excel COM-interop script
# Applies numeric corrections to the profit-calculation columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 1237.5
$ws.Range("J12").Value = 1225
$ws.Range("L12").Value = 1225
$ws.Range("N12").Value = -1565

# Row 41
$ws.Range("H41").Value = 531.44446
$ws.Range("I41").Value = 397
$ws.Range("J41").Value = 699.5
$ws.Range("K41").Value = 397
$ws.Range("L41").Value = 699.5
$ws.Range("M41").Value = 43
$ws.Range("N41").Value = -1579.5

# Row 107
$ws.Range("H107").Value = 429.55
$ws.Range("I107").Value = 352.41177
$ws.Range("J107").Value = 866.6667
$ws.Range("K107").Value = 352.41177
$ws.Range("L107").Value = 866.6667
$ws.Range("M107").Value = 1567.58823
$ws.Range("N107").Value = -4706.6667

# Row 112
$ws.Range("H112").Value = 3168.45
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3168.45
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 9505.349999999999
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -11721.35

# Row 129
$ws.Range("H129").Value = 917.86884
$ws.Range("I129").Value = 883.7
$ws.Range("J129").Value = 924.5685999999999
$ws.Range("K129").Value = 2651.1
$ws.Range("L129").Value = 2773.7058
$ws.Range("M129").Value = 2348.9
$ws.Range("N129").Value = -12773.7058

# Row 135
$ws.Range("H135").Value = 422.04544
$ws.Range("I135").Value = 349.22223
$ws.Range("K135").Value = 3143.00007
$ws.Range("M135").Value = -608.0000700000001

# Row 138
$ws.Range("H138").Value = 3687.9714
$ws.Range("I138").Value = 3118.64
$ws.Range("J138").Value = 5111.3
$ws.Range("K138").Value = 9355.92
$ws.Range("L138").Value = 15333.9
$ws.Range("M138").Value = -4215.92
$ws.Range("N138").Value = -25613.9

# Row 140
$ws.Range("H140").Value = 77699.25
$ws.Range("J140").Value = 77699.25
$ws.Range("L140").Value = 77699.25
$ws.Range("N140").Value = -88059.25

# Row 141
$ws.Range("H141").Value = 1870040.1
$ws.Range("I141").Value = 3503049
$ws.Range("J141").Value = 3744.1428
$ws.Range("K141").Value = 10509147
$ws.Range("L141").Value = 11232.4284
$ws.Range("M141").Value = -10503967
$ws.Range("N141").Value = -21592.4284

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4197.035
$ws.Range("I32").Value = 3332.0417
$ws.Range("K32").Value = 3332.0417
$ws.Range("M32").Value = -3045.0417

# Row 74
$ws.Range("H74").Value = 1000
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 1000
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -126
$ws.Range("N74").Value = -2748

# Row 77
$ws.Range("H77").Value = 1000
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 5000
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -632
$ws.Range("N77").Value = -13736

# Row 102
$ws.Range("H102").Value = 1165.2222
$ws.Range("J102").Value = 1000
$ws.Range("L102").Value = 1000
$ws.Range("N102").Value = -4244

# Row 110
$ws.Range("H110").Value = 1122.8572
$ws.Range("I110").Value = 976.6667
$ws.Range("K110").Value = 976.6667
$ws.Range("M110").Value = 1068.3333

# Row 122
$ws.Range("H122").Value = 1323.1428
$ws.Range("I122").Value = 1358.6957
$ws.Range("J122").Value = 1159.6
$ws.Range("K122").Value = 4076.0871
$ws.Range("L122").Value = 3478.8
$ws.Range("M122").Value = -1626.0871
$ws.Range("N122").Value = -8378.799999999999

# Row 123
$ws.Range("H123").Value = 62497.5
$ws.Range("J123").Value = 62497.5
$ws.Range("L123").Value = 62497.5
$ws.Range("N123").Value = -72297.5

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 5755.436
$ws.Range("I134").Value = 6326.839
$ws.Range("K134").Value = 18980.517
$ws.Range("M134").Value = -16445.517

$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 1000
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1340

# Row 22
$ws.Range("H22").Value = 1378.75
$ws.Range("I22").Value = 1082.5
$ws.Range("J22").Value = 1675
$ws.Range("K22").Value = 1082.5
$ws.Range("L22").Value = 1675
$ws.Range("M22").Value = -732.5
$ws.Range("N22").Value = -2375

# Row 31
$ws.Range("H31").Value = 1851.6522
$ws.Range("I31").Value = 1017.16
$ws.Range("J31").Value = 2845.0952
$ws.Range("K31").Value = 1017.16
$ws.Range("L31").Value = 2845.0952
$ws.Range("M31").Value = -722.16
$ws.Range("N31").Value = -3435.0952

# Row 32
$ws.Range("H32").Value = 12641.2
$ws.Range("I32").Value = 12641.2
$ws.Range("K32").Value = 12641.2
$ws.Range("M32").Value = -12325.2

# Row 34
$ws.Range("H34").Value = 1851.6522
$ws.Range("I34").Value = 1017.16
$ws.Range("J34").Value = 2845.0952
$ws.Range("K34").Value = 1017.16
$ws.Range("L34").Value = 2845.0952
$ws.Range("M34").Value = -815.16
$ws.Range("N34").Value = -3249.0952

# Row 62
$ws.Range("H62").Value = 3123.3076
$ws.Range("I62").Value = 3099.6667
$ws.Range("K62").Value = 3099.6667
$ws.Range("M62").Value = -2475.6667

# Row 65
$ws.Range("H65").Value = 3123.3076
$ws.Range("I65").Value = 3099.6667
$ws.Range("K65").Value = 15498.3335
$ws.Range("M65").Value = -12378.3335

# Row 107
$ws.Range("H107").Value = 526.2778
$ws.Range("I107").Value = 388.9091
$ws.Range("J107").Value = 742.1429000000001
$ws.Range("K107").Value = 388.9091
$ws.Range("L107").Value = 742.1429000000001
$ws.Range("M107").Value = 1531.0909
$ws.Range("N107").Value = -4582.1429

# Row 122
$ws.Range("H122").Value = 1957.3334
$ws.Range("I122").Value = 1276
$ws.Range("J122").Value = 2502.4
$ws.Range("K122").Value = 3828
$ws.Range("L122").Value = 7507.200000000001
$ws.Range("M122").Value = -1378
$ws.Range("N122").Value = -12407.2

# Row 134
$ws.Range("H134").Value = 2505.6086
$ws.Range("I134").Value = 1727.5
$ws.Range("K134").Value = 5182.5
$ws.Range("M134").Value = -2647.5

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 999
$ws.Range("J2").Value = 998.8333
$ws.Range("L2").Value = 5992.9998
$ws.Range("N2").Value = -6218.9998

# Row 12
$ws.Range("H12").Value = 62.666668
$ws.Range("I12").Value = 22.5
$ws.Range("J12").Value = 108.57143
$ws.Range("K12").Value = 67.5
$ws.Range("L12").Value = 325.71429
$ws.Range("M12").Value = 105.5
$ws.Range("N12").Value = -671.71429

# Row 17
$ws.Range("H17").Value = 10000
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 30000
$ws.Range("N17").Value = -30338

# Row 39
$ws.Range("H39").Value = 2542.2856
$ws.Range("J39").Value = 3045
$ws.Range("L39").Value = 9135
$ws.Range("N39").Value = -9723

# Row 50
$ws.Range("H50").Value = 71479330
$ws.Range("I50").Value = 116580.164
$ws.Range("J50").Value = 125001390
$ws.Range("K50").Value = 349740.492
$ws.Range("L50").Value = 375004170
$ws.Range("M50").Value = -349259.492
$ws.Range("N50").Value = -375005132

# Row 53
$ws.Range("H53").Value = 71479330
$ws.Range("I53").Value = 116580.164
$ws.Range("J53").Value = 125001390
$ws.Range("K53").Value = 349740.492
$ws.Range("L53").Value = 375004170
$ws.Range("M53").Value = -349259.492
$ws.Range("N53").Value = -375005132

# Row 55
$ws.Range("H55").Value = 100004
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# Row 107
$ws.Range("H107").Value = 1868.5625
$ws.Range("J107").Value = 1879.1666
$ws.Range("L107").Value = 5637.4998
$ws.Range("N107").Value = -9477.4998

# Row 118
$ws.Range("H118").Value = 55558870
$ws.Range("I118").Value = 166667620
$ws.Range("J118").Value = 4500
$ws.Range("K118").Value = 500002860
$ws.Range("L118").Value = 13500
$ws.Range("M118").Value = -500001617
$ws.Range("N118").Value = -15986

# Row 131
$ws.Range("H131").Value = 14306909
$ws.Range("J131").Value = 22448.516
$ws.Range("L131").Value = 67345.548
$ws.Range("N131").Value = -77425.548

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 36833
$ws.Range("I10").Value = 5249.5
$ws.Range("K10").Value = 5249.5
$ws.Range("M10").Value = -5080.5

# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1937
$ws.Range("I22").Value = 1749
$ws.Range("J22").Value = 1999.6666
$ws.Range("K22").Value = 1749
$ws.Range("L22").Value = 1999.6666
$ws.Range("M22").Value = -1454
$ws.Range("N22").Value = -2589.6666

# Row 27
$ws.Range("H27").Value = 1937
$ws.Range("I27").Value = 1749
$ws.Range("J27").Value = 1999.6666
$ws.Range("K27").Value = 1749
$ws.Range("L27").Value = 1999.6666
$ws.Range("M27").Value = -1642
$ws.Range("N27").Value = -2213.6666

# Row 132
$ws.Range("H132").Value = 3565.2222
$ws.Range("I132").Value = 2167.6667
$ws.Range("J132").Value = 3739.9167
$ws.Range("K132").Value = 6503.000100000001
$ws.Range("L132").Value = 11219.7501
$ws.Range("M132").Value = -3973.000100000001
$ws.Range("N132").Value = -16279.7501

# Row 136
$ws.Range("H136").Value = 5097.136
$ws.Range("I136").Value = 3566.5833
$ws.Range("K136").Value = 10699.7499
$ws.Range("M136").Value = -8149.749899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2374
$ws.Range("I132").Value = 998.5
$ws.Range("J132").Value = 3749.5
$ws.Range("K132").Value = 2995.5
$ws.Range("L132").Value = 11248.5
$ws.Range("M132").Value = -465.5
$ws.Range("N132").Value = -16308.5

# Row 136
$ws.Range("H136").Value = 11339843
$ws.Range("I136").Value = 26457360
$ws.Range("J136").Value = 1706.6072
$ws.Range("K136").Value = 79372080
$ws.Range("L136").Value = 5119.821599999999
$ws.Range("M136").Value = -79369530
$ws.Range("N136").Value = -10219.8216
